$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes what used to be row 3's data (Abhay Singh, Carrom Board, Faculty, no F2 value)
$ws.Range("A2").Value = "Abhay Singh"
$ws.Range("B2").Value = "abhaysinghktp800@gmail.com"
$ws.Range("C2").Value = "Carrom Board"
$ws.Range("D2").Value = "'09588014420"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "Faculty"
$ws.Range("F2").ClearContents()

# Row 3 becomes what used to be row 2's data (Aman choudhary, Badminton, Student, BCA)
$ws.Range("A3").Value = "Aman choudhary"
$ws.Range("B3").Value = "gauravbadaliya420@gmail.com"
$ws.Range("C3").Value = "Badminton"
$ws.Range("D3").Value = "'07988476183"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "Student"
$ws.Range("F3").Value = "BCA"

# Row 4: F4 changes from BCA to BTECH
$ws.Range("F4").Value = "BTECH"
